$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "treatment"
$ws.Range("B1").Value = "bone_mineral_density"
$ws.Range("C1").Value = "description"

$ws.Range("A2").Select() | Out-Null
